# Updated symbol list on Sat Dec 24 03:47:46 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value as exact text, without leaving a
# "number stored as text" style behind (keeps cell's original style).
function Set-ExactText {
    param($cell, [string]$text)
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

# --- Price (column D) updates -------------------------------------------
Set-ExactText "D2"  "245.90"
Set-ExactText "D4"  "5.357"
Set-ExactText "D5"  "0.05925"
Set-ExactText "D6"  "3.394"
Set-ExactText "D8"  "0.8131"
Set-ExactText "D9"  "0.9598"
Set-ExactText "D10" "0.1429"
Set-ExactText "D11" "0.07402"
Set-ExactText "D12" "0.03468"
Set-ExactText "D13" "0.03042"
Set-ExactText "D16" "0.001588"
Set-ExactText "D17" "0.04813"
Set-ExactText "D19" "0.006032"
Set-ExactText "D20" "0.004083"
Set-ExactText "D21" "0.0009908"
Set-ExactText "D22" "0.00009703"
Set-ExactText "D23" "3.743"
Set-ExactText "D40" "0.03940"
Set-ExactText "D41" "0.006432"
Set-ExactText "D42" "0.1074"
Set-ExactText "D43" "0.002711"
Set-ExactText "D44" "0.005899"
Set-ExactText "D45" "0.00005303"
Set-ExactText "D48" "0.04661"

# --- Rows 14 and 15 swap places (MCDex <-> BitMartToken) ----------------
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-ExactText "D14" "0.09404"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-ExactText "D15" "4.001"
$ws.Range("E15").Value = "14MCDexMCB"

# --- "Bestin24h" suffix moved from row 41 to row 47 ----------------------
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
